$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")

$ws.Range("A2").Value = "Lavern"
$ws.Range("B2").Value = "Wiza"
$ws.Range("C2").Value = "erasmo.muller@yahoo.com"
$ws.Range("D2").Value = "kqmp15ef"
